# Auto-generated Excel COM-interop script to apply crypto price/volume updates
# (mirrors the Mon Sep 18 22:50:31 UTC 2023 GitHub Actions refresh of cryptos.xlsx)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.881.33'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.72%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.642.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.74%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.56'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.88%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0623'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.39%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.80'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +4.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0847'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.872.04'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.64%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.658.10'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.00%  '
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.529'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.31'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.67%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.879.85'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.92%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '219.50'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +3.92%  '
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.72%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.63'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.21%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.19'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.32%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.04'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.41'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +5.92%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.82'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +1.92%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.08%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.00'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.56'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.53%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.246.05'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.538'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.62%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.833'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.89%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.806'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.29%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.783.17'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.75%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '60.77'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.81%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0514'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.20%  '
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0974'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.98%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.62'
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.405'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.22%  '
